# Daily attendance processing - 2026-01-15 11:35:17
# Normalizes the "Recorded By" column (G) so that the canonical recorder
# (an email address) is listed last, after any system-generated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -lt 2) {
        continue
    }

    if ($parts[0] -like "*@*") {
        $rest = $parts[1..($parts.Count - 1)]
        $newParts = $rest + $parts[0]
        $cell.Value = $newParts -join ", "
    }
}
